$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flat list: each logical row contributes 8 numbers -> RowNumber, then values for
# columns B, C, D, E, F, I, N (the only columns whose values change in this edit).
$cols = @("B","C","D","E","F","I","N")
$stride = $cols.Length + 1

$data = @(
    2,1.43180100408847,0.5561946450787332,0.01898470219389026,0.418868025178611,1.22420842577543,0.6851918156816694,0.9873592261007289,
    3,1.270843218843311,0.4868603933235249,0.01914739503100904,0.3649584804624908,1.168444706594499,0.6700495542047094,1.003888397040395,
    4,1.172501167612211,0.4444327630589555,0.01926975491615224,0.3320111203335046,1.135248235746545,0.6613475268064306,1.014555381492585,
    5,1.132544577260489,0.4271769619589918,0.0193251675234194,0.3186199087233206,1.121978765897921,0.6579493326846304,1.01903186804314,
    6,1.125916869622642,0.4243136222090698,0.01933470107326229,0.3163983286809184,1.119790873934249,0.6573939500089239,1.019782994033985,
    7,1.171961823827075,0.4441999115416593,0.01927047988786512,0.331830384154614,1.135068237959004,0.661301100826023,1.014615229045823,
    8,1.376199399961763,0.5322568711487179,0.01903607571363963,0.4002460886429873,1.204762350224854,0.6798461201412564,0.9929505784159041,
    9,1.780764153406892,0.7061947151043455,0.01875920105149831,0.535788821193421,1.349887909244501,0.7210197532924383,0.9546052659177207,
    10,2.08077230338705,0.8349296160035919,0.01867358378142114,0.636450542645548,1.461933247349464,0.7543200522087972,0.9289952987538861,
    11,2.217924869523586,0.8937386091240569,0.01866162578351194,0.6825312684266436,1.514142216851695,0.7701580401969039,0.9179096619487979,
    12,2.269963472863196,0.9160465645317117,0.0186610976509769,0.7000266628881207,1.534095255452371,0.7762567557253419,0.9137937574140693,
    13,2.2587514234379,0.9112404005356325,0.0186610315642568,0.6962566223281641,1.529789816759688,0.7749387535038181,0.9146765315895671,
    14,2.2222040482431,0.8955731150404063,0.01866150147474244,0.6839696872497427,1.515780079165381,0.7706577447913077,0.9175693961619196,
    15,2.199831132061433,0.8859815169456056,0.01866231381580974,0.6764496476126993,1.507222628932141,0.7680487461037444,0.9193520606761254,
    16,2.071823072738027,0.8310915125155702,0.01867491848107505,0.633445221417773,1.458546580191637,0.7532990526687229,0.929731199063383,
    17,1.993470925925237,0.7974834466871812,0.01868964761732528,0.6071403153444521,1.429005889067014,0.7444286503622948,0.9362436128909799,
    18,1.948468322837186,0.7781760405782165,0.01870065255525333,0.592037403652057,1.412131211511991,0.7393914045231753,0.9400424614798908,
    19,1.933241982149468,0.7716427625505276,0.01870481012150549,0.5869283286249072,1.40643758340525,0.737696942551608,0.9413377783360701,
    20,2.001805043607987,0.8010586757072247,0.0186878166532054,0.6099376967771946,1.432138472512435,0.7453662010517519,0.935544855845535,
    21,2.232936096178037,0.9001739156577742,0.01866125392796647,0.6875773871021522,1.519890085992898,0.7719124173112704,0.9167174600614842,
    22,2.38458974426311,0.9651758285419874,0.01866726640396266,0.7385872800035287,1.578307509139734,0.7898528354055259,0.9048909335181037,
    23,2.303593272698947,0.9304616225160771,0.01866187882006187,0.7113364733707357,1.547029931490101,0.7802229355520041,0.9111589380551806,
    24,1.998037055126417,0.7994422702759607,0.01868863653876218,0.6086729382376603,1.430721893424504,0.7449421399710729,0.9358605933551623,
    25,1.670853665900438,0.6589872455524528,0.01881388412277829,0.4989489557592321,1.309695057901877,0.7093543365496231,0.9645315964188583
)

for ($j = 0; $j -lt $data.Length; $j += $stride) {
    $row = $data[$j]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $data[$j + 1 + $i]
        $addr = "$col$row"
        $ws.Range($addr).Value = $value
    }
}
